$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we are about to touch to Text format
# so that numeric-looking strings (e.g. "6.80", "1.00") are preserved as
# literal text instead of being auto-converted to numbers by Excel.
$dCells = @("D2","D3","D5","D6","D8","D9","D12","D13","D15","D16","D18","D19","D20","D21","D22","D23","D24","D25","D26","D28","D31","D32","D34","D35","D36","D37","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '67.201.83'
$ws.Range("E2").Value = '  +4.98%  '

# Row 3
$ws.Range("D3").Value = '3.240.20'
$ws.Range("E3").Value = '  +2.97%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").Value = '576.55'
$ws.Range("E5").Value = '  +2.77%  '

# Row 6
$ws.Range("D6").Value = '178.44'
$ws.Range("E6").Value = '  +6.86%  '

# Row 8
$ws.Range("D8").Value = '0.598'
$ws.Range("E8").Value = '  -3.45%  '

# Row 9
$ws.Range("D9").Value = '3.238.05'
$ws.Range("E9").Value = '  +2.97%  '

# Row 10
$ws.Range("E10").Value = '  +5.35%  '

# Row 11
$ws.Range("E11").Value = '  +3.11%  '

# Row 12
$ws.Range("D12").Value = '0.412'
$ws.Range("E12").Value = '  +4.98%  '

# Row 13
$ws.Range("D13").Value = '3.802.69'
$ws.Range("E13").Value = '  +2.99%  '

# Row 14
$ws.Range("E14").Value = '  +0.95%  '

# Row 15
$ws.Range("D15").Value = '27.83'
$ws.Range("E15").Value = '  +3.11%  '

# Row 16
$ws.Range("D16").Value = '67.140.11'
$ws.Range("E16").Value = '  +4.92%  '

# Row 17
$ws.Range("E17").Value = '  +3.38%  '

# Row 18
$ws.Range("D18").Value = '3.242.33'
$ws.Range("E18").Value = '  +2.89%  '

# Row 19
$ws.Range("D19").Value = '5.81'
$ws.Range("E19").Value = '  +2.03%  '

# Row 20
$ws.Range("D20").Value = '13.33'
$ws.Range("E20").Value = '  +3.95%  '

# Row 21
$ws.Range("D21").Value = '372.95'
$ws.Range("E21").Value = '  +6.47%  '

# Row 22
$ws.Range("D22").Value = '7.55'
$ws.Range("E22").Value = '  +5.91%  '

# Row 23
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.13%  '

# Row 24
$ws.Range("D24").Value = '71.23'
$ws.Range("E24").Value = '  +4.79%  '

# Row 25
$ws.Range("D25").Value = '0.508'
$ws.Range("E25").Value = '  +1.88%  '

# Row 26
$ws.Range("D26").Value = '3.384.28'
$ws.Range("E26").Value = '  +2.41%  '

# Row 27
$ws.Range("E27").Value = '  +3.74%  '

# Row 28
$ws.Range("D28").Value = '9.65'
$ws.Range("E28").Value = '  +1.20%  '

# Row 29
$ws.Range("E29").Value = '  +2.98%  '

# Row 31
$ws.Range("D31").Value = '1.97'
$ws.Range("E31").Value = '  +4.81%  '

# Row 32
$ws.Range("D32").Value = '5.61'
$ws.Range("E32").Value = '  +2.96%  '

# Row 33
$ws.Range("E33").Value = '  +3.42%  '

# Row 34
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.01%  '

# Row 35
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").Value = '1.27'
$ws.Range("E35").Value = '  +7.60%  '

# Row 36
$ws.Range("D36").Value = '6.80'
$ws.Range("E36").Value = '  +3.54%  '

# Row 37
$ws.Range("D37").Value = '163.70'
$ws.Range("E37").Value = '  +6.66%  '

# Row 38
$ws.Range("E38").Value = '  +4.74%  '

# Row 39
$ws.Range("D39").Value = '0.859'
$ws.Range("E39").Value = '  +5.72%  '

# Row 40
$ws.Range("D40").Value = '1.84'
$ws.Range("E40").Value = '  +9.59%  '

# Row 41
$ws.Range("D41").Value = '6.78'
$ws.Range("E41").Value = '  +13.64%  '

# Row 42
$ws.Range("D42").Value = '26.65'
$ws.Range("E42").Value = '  +2.67%  '

# Row 43
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '364.34'
$ws.Range("E43").Value = '  +16.15%  '

# Row 44
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '2.60'
$ws.Range("E44").Value = '  +6.07%  '

# Row 45
$ws.Range("D45").Value = '2.716.30'
$ws.Range("E45").Value = '  +5.07%  '

# Row 46
$ws.Range("D46").Value = '4.39'
$ws.Range("E46").Value = '  +5.82%  '

# Row 47
$ws.Range("D47").Value = '25.81'
$ws.Range("E47").Value = '  +10.01%  '

# Row 48
$ws.Range("D48").Value = '40.34'
$ws.Range("E48").Value = '  +2.90%  '

# Row 49
$ws.Range("D49").Value = '0.0673'
$ws.Range("E49").Value = '  +5.02%  '

# Row 50
$ws.Range("D50").Value = '0.0277'
$ws.Range("E50").Value = '  +3.69%  '

# Row 51
$ws.Range("E51").Value = '  +1.05%  '

# Clear the temporary text-format override so the cells end up with
# no explicit style, matching their original style-less state.
foreach ($addr in $dCells) {
    $ws.Range($addr).ClearFormats()
}
